$d = $word.ActiveDocument

# 1. Replace the ID placeholder text (also merges/cleans the trailing space run)
$d.Content.Find.Execute("**ID__AFFARS_pgi_5301_topic_50__ID** ", $false, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_USAFA_PGI_5301__ID**", 2)

$p = $d.Paragraphs(1)

# 2. Update left indent from 120 twips (6pt) to 225 twips (11.25pt)
$p.Format.LeftIndent = 11.25

# 3. Add a paragraph border (top/left/bottom/right) with 5-twip spacing, no line
$b = $p.Range.Borders
$b.DistanceFromTop = 5
$b.DistanceFromBottom = 5
$b.DistanceFromLeft = 5
$b.DistanceFromRight = 5
